$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$lo = $ws.ListObjects.Item("Table1")

# Fix the typo'd header ("No. of (Offenses" -> "No. of (Burglary)")
$ws.Range("G1").Value = "No. of (Burglary)"

# Add a new table column for non-violent crime counts
$lo.ListColumns.Add() | Out-Null
$ws.Range("H1").Value = "No. of Non-Violent Crimes"

# Fill in ORI Codes (column E) for the states I looked up
$ws.Range("E31").Value = "NJ0111100"
$ws.Range("E32").Value = "NM0260100"
$ws.Range("E33").Value = "NY0010100"
$ws.Range("E34").Value = "NC0920100"
$ws.Range("E35").Value = "ND0080100"
$ws.Range("E36").Value = "OHCOP0000"
$ws.Range("E37").Value = "OK0550400"
$ws.Range("E38").Value = "OR0240200"
$ws.Range("E39").Value = "PA0220200"

# Match the font used by the pasted-in ORI codes
$ori = $ws.Range("E31:E40")
$ori.Font.Name = "Arial Unicode MS"
$ori.Font.Size = 10
$ori.Font.Color = 0

# Row 37 (Oklahoma) got manually resized while the data was being entered
$ws.Rows.Item(37).RowHeight = 24

# Widen the burglary/non-violent-crime columns so the longer headers fit
$ws.Columns.Item(7).ColumnWidth = 26.85546875
$ws.Columns.Item(8).ColumnWidth = 29.28515625

$ws.Range("F44").Select() | Out-Null
